$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 53 (01-01-2021) ---
$ws.Range("B53").Value = 86.09999999999999
$ws.Range("D53").Value = 116.8

# --- Add new row 54 (01-04-2021) ---
# Writing a date-like string such as "01-04-2021" directly into .Value/.Formula
# gets auto-recognized as a date by Excel and converted to a serial number.
# To keep it as plain text (matching the other "Serie" cells in column A which
# are stored as shared strings), build the text via a formula (formula results
# that are text are never re-interpreted as dates) in a scratch cell, then
# copy only the resulting value into A54, then clear the scratch cell.
$scratch = $ws.Range("ZZ1")
$scratch.Formula = '=""&"01-04-2021"'
$scratch.Copy()
$ws.Range("A54").PasteSpecial(-4163)
$scratch.ClearContents()

$ws.Range("B54").Value = 85.5
$ws.Range("C54").Value = 30.7
$ws.Range("D54").Value = 116.2

$excel.CutCopyMode = 0
